$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column AD (30) - header + weekly data for 25_05_2021
$ws.Cells.Item(1, 30).Value = "25_05_2021"

$ws.Cells.Item(2, 30).Value = 12
$ws.Cells.Item(3, 30).Value = 17
$ws.Cells.Item(4, 30).Value = 37
$ws.Cells.Item(5, 30).Value = 48
$ws.Cells.Item(6, 30).Value = 124
$ws.Cells.Item(7, 30).Value = 287
$ws.Cells.Item(8, 30).Value = 417
$ws.Cells.Item(9, 30).Value = 574
$ws.Cells.Item(10, 30).Value = 178
$ws.Cells.Item(11, 30).Value = 13

$ws.Range("AD12").Formula = "=SUM(AD2:AD11)"

# Update the selected cell / view as recorded in the workbook
$ws.Range("AD15").Select()
